$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("C1").Value = "Precio"
$ws.Range("D1").Value = "Cantidad"

$ws.Range("B2").Value = "Macbook air M3"
$ws.Range("B3").Value = "Macbook Pro M4"

$ws.Range("C3").Value = "2000$"
$ws.Range("C2").Value = "1000$"

$ws.Range("D2").Value = "10"
$ws.Range("D3").Value = "20"

$ws.Range("D6").Select()
$excel.ActiveWindow.Zoom = 112

$ws.Columns("A").ColumnWidth = 12.2
$ws.Columns("B").ColumnWidth = 15.6
